# Updates the cryptos price/volume table to the latest scrape.
# Rows 2-46: refresh Price (D) and Volume(1h) (E) values in place.
# Rows 47-51: the lowest-ranked coin (EnergySwap) drops off the list,
# the remaining four coins shift up one row, and a new coin
# (ApeXProtocol) is appended at the bottom with its own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    # Force the cell to stay a text value (many of these look like
    # numbers, e.g. "1.00" or "0.586") without leaving a residual
    # explicit cell style behind once the custom number format is
    # no longer needed.
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

# --- Rows 2-46: price & volume refresh ---

Set-TextValue $ws.Range("D2") '63.764.84'
Set-TextValue $ws.Range("E2") '  -5.59%  '

Set-TextValue $ws.Range("D3") '3.288.34'
Set-TextValue $ws.Range("E3") '  -7.93%  '

Set-TextValue $ws.Range("E4") '  +0.11%  '

Set-TextValue $ws.Range("D5") '180.35'
Set-TextValue $ws.Range("E5") '  -9.75%  '

Set-TextValue $ws.Range("D6") '514.75'
Set-TextValue $ws.Range("E6") '  -9.31%  '

Set-TextValue $ws.Range("D7") '0.586'
Set-TextValue $ws.Range("E7") '  -4.27%  '

Set-TextValue $ws.Range("D8") '3.279.04'
Set-TextValue $ws.Range("E8") '  -8.01%  '

Set-TextValue $ws.Range("E9") '  +0.06%  '

Set-TextValue $ws.Range("D10") '0.610'
Set-TextValue $ws.Range("E10") '  -9.37%  '

Set-TextValue $ws.Range("D11") '57.44'
Set-TextValue $ws.Range("E11") '  -4.21%  '

Set-TextValue $ws.Range("E12") '  -11.65%  '

Set-TextValue $ws.Range("D13") '0.0000248'
Set-TextValue $ws.Range("E13") '  -11.13%  '

Set-TextValue $ws.Range("D14") '8.96'
Set-TextValue $ws.Range("E14") '  -12.40%  '

Set-TextValue $ws.Range("D15") '3.833.73'
Set-TextValue $ws.Range("E15") '  -8.05%  '

Set-TextValue $ws.Range("E16") '  -4.61%  '

Set-TextValue $ws.Range("D17") '3.314.26'
Set-TextValue $ws.Range("E17") '  -7.67%  '

Set-TextValue $ws.Range("D18") '63.564.68'
Set-TextValue $ws.Range("E18") '  -5.75%  '

Set-TextValue $ws.Range("D19") '17.00'
Set-TextValue $ws.Range("E19") '  -10.58%  '

Set-TextValue $ws.Range("D20") '10.73'
Set-TextValue $ws.Range("E20") '  -11.83%  '

Set-TextValue $ws.Range("D21") '0.937'
Set-TextValue $ws.Range("E21") '  -11.16%  '

Set-TextValue $ws.Range("D22") '366.19'
Set-TextValue $ws.Range("E22") '  -8.73%  '

Set-TextValue $ws.Range("D23") '79.37'
Set-TextValue $ws.Range("E23") '  -6.07%  '

Set-TextValue $ws.Range("D24") '3.63'
Set-TextValue $ws.Range("E24") '  -12.69%  '

Set-TextValue $ws.Range("D25") '10.59'
Set-TextValue $ws.Range("E25") '  -17.38%  '

Set-TextValue $ws.Range("D26") '6.00'
Set-TextValue $ws.Range("E26") '  -1.71%  '

Set-TextValue $ws.Range("D27") '3.69'
Set-TextValue $ws.Range("E27") '  -5.07%  '

Set-TextValue $ws.Range("D28") '2.60'
Set-TextValue $ws.Range("E28") '  -9.58%  '

Set-TextValue $ws.Range("D29") '11.04'
Set-TextValue $ws.Range("E29") '  -10.61%  '

Set-TextValue $ws.Range("D30") '8.21'
Set-TextValue $ws.Range("E30") '  -10.46%  '

Set-TextValue $ws.Range("D31") '642.84'
Set-TextValue $ws.Range("E31") '  -4.16%  '

Set-TextValue $ws.Range("D32") '28.18'
Set-TextValue $ws.Range("E32") '  -9.99%  '

Set-TextValue $ws.Range("D33") '6.59'
Set-TextValue $ws.Range("E33") '  -13.96%  '

Set-TextValue $ws.Range("D34") '10.94'
Set-TextValue $ws.Range("E34") '  -9.26%  '

Set-TextValue $ws.Range("D35") '58.94'
Set-TextValue $ws.Range("E35") '  -6.95%  '

Set-TextValue $ws.Range("E36") '  -9.49%  '

Set-TextValue $ws.Range("D37") '1.00'
Set-TextValue $ws.Range("E37") '  +0.01%  '

Set-TextValue $ws.Range("D38") '35.49'
Set-TextValue $ws.Range("E38") '  -14.00%  '

Set-TextValue $ws.Range("D39") '0.369'
Set-TextValue $ws.Range("E39") '  -9.26%  '

Set-TextValue $ws.Range("D40") '0.999'
Set-TextValue $ws.Range("E40") '  -0.07%  '

Set-TextValue $ws.Range("D41") '0.121'
Set-TextValue $ws.Range("E41") '  -8.68%  '

Set-TextValue $ws.Range("D42") '2.771.16'
Set-TextValue $ws.Range("E42") '  -12.66%  '

Set-TextValue $ws.Range("D43") '2.63'
Set-TextValue $ws.Range("E43") '  -17.13%  '

Set-TextValue $ws.Range("D44") '2.59'
Set-TextValue $ws.Range("E44") '  -7.43%  '

Set-TextValue $ws.Range("D45") '0.0₃0606'
Set-TextValue $ws.Range("E45") '  -19.70%  '

Set-TextValue $ws.Range("D46") '0.0381'
Set-TextValue $ws.Range("E46") '  -6.85%  '

# --- Rows 47-51: list shift + new entry ---

$ws.Range("B47").Value = 'Fetch.AI'
$ws.Range("C47").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range("D47") '2.26'
Set-TextValue $ws.Range("E47") '  -16.03%  '

$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range("D48") '0.122'
Set-TextValue $ws.Range("E48") '  -6.13%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range("D49") '133.81'
Set-TextValue $ws.Range("E49") '  -3.39%  '

$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D50") '2.62'
Set-TextValue $ws.Range("E50") '  -2.78%  '

$ws.Range("B51").Value = 'ApeXProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws.Range("D51") '2.76'
Set-TextValue $ws.Range("E51") '  -10.89%  '

